$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.09713344452125153
$ws.Range("E2").Value = 0.1649182917956413
$ws.Range("H2").Value = -0.07287325603580108
$ws.Range("J2").Value = 4.664256048167697
$ws.Range("K2").Value = 1.255846506602602

# Row 3
$ws.Range("C3").Value = 0.1357309662849703
$ws.Range("E3").Value = 0.4304909107745751
$ws.Range("H3").Value = -0.1291799776589777

# Row 4
$ws.Range("C4").Value = 0.03602140854847032
$ws.Range("E4").Value = 0.715047776010967
$ws.Range("H4").Value = 0.5157574252241206

# Row 5
$ws.Range("C5").Value = 0.1673970425776976
$ws.Range("E5").Value = 1.040992192653427

# Row 6
$ws.Range("C6").Value = 0.1439811213567938
$ws.Range("E6").Value = 1.363918250517865
$ws.Range("H6").Value = -0.1368514987685632

# Row 7
$ws.Range("C7").Value = -0.1291799776589777
$ws.Range("E7").Value = 1.166125193703865

# Row 8
$ws.Range("C8").Value = 0.04662796151971213
$ws.Range("E8").Value = 0.1100062521651354
$ws.Range("H8").Value = 0.3420446477696094

# Row 9
$ws.Range("C9").Value = 0.2280148978263631
$ws.Range("E9").Value = 0.4903809066616679

# Row 10
$ws.Range("C10").Value = -0.1368514987685632
$ws.Range("E10").Value = 1.210433602295166
$ws.Range("H10").Value = 0.2034184511261679

# Row 11
$ws.Range("C11").Value = -0.07287325603580108
$ws.Range("E11").Value = 0.9772874142478776

# Row 12
$ws.Range("C12").Value = 0.2063136814846392
$ws.Range("E12").Value = 0.1751184849636601
$ws.Range("H12").Value = -0.09933819709753919

# Row 13
$ws.Range("C13").Value = -0.09933819709753919
$ws.Range("E13").Value = 0.4930194929127935
